# Update the "featuredNewsPage" test data: the 4th featured-news headline
# (row 5, column B of Sheet1) is replaced with a new headline.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "Cognizant shines at the Times Group Global Business Summit (GBS)"
